$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for the Price (D) and Volume(1h) (E) columns so
# numeric-looking strings (e.g. "1.00", "67.20") are not coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '42.088.05'
$ws.Range('E2').Value = '  -0.90%  '
$ws.Range('D3').Value = '2.250.15'
$ws.Range('E3').Value = '  -1.46%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '306.44'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('D6').Value = '96.84'
$ws.Range('E6').Value = '  -1.40%  '
$ws.Range('D7').Value = '0.524'
$ws.Range('E7').Value = '  -1.37%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '0.488'
$ws.Range('E9').Value = '  -1.43%  '
$ws.Range('D10').Value = '34.89'
$ws.Range('E10').Value = '  -3.43%  '
$ws.Range('D11').Value = '0.0811'
$ws.Range('E11').Value = '  +1.42%  '
$ws.Range('E12').Value = '  +1.18%  '
$ws.Range('D13').Value = '6.78'
$ws.Range('E13').Value = '  +0.94%  '
$ws.Range('D14').Value = '2.599.26'
$ws.Range('E14').Value = '  -1.38%  '
$ws.Range('D15').Value = '14.43'
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('D16').Value = '2.245.41'
$ws.Range('E16').Value = '  -2.03%  '
$ws.Range('D17').Value = '0.779'
$ws.Range('E17').Value = '  -2.55%  '
$ws.Range('D18').Value = '41.970.40'
$ws.Range('E18').Value = '  -0.91%  '
$ws.Range('D19').Value = '12.16'
$ws.Range('E19').Value = '  -3.52%  '
$ws.Range('D20').Value = '0.0₃0902'
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('D21').Value = '5.92'
$ws.Range('D22').Value = '67.20'
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('D23').Value = '235.38'
$ws.Range('E23').Value = '  -2.44%  '
$ws.Range('E24').Value = '  -1.77%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = '1.02'
$ws.Range('E25').Value = '  +1.59%  '
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').Value = '1.94'
$ws.Range('E26').Value = '  -0.82%  '
$ws.Range('B27').Value = 'InjectiveProtocol'
$ws.Range('C27').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D27').Value = '37.95'
$ws.Range('E27').Value = '  +0.38%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '23.29'
$ws.Range('E28').Value = '  -2.69%  '
$ws.Range('D29').Value = '2.12'
$ws.Range('E29').Value = '  +0.67%  '
$ws.Range('D30').Value = '9.47'
$ws.Range('E30').Value = '  -0.87%  '
$ws.Range('D31').Value = '167.11'
$ws.Range('E31').Value = '  +4.63%  '
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('D33').Value = '5.16'
$ws.Range('E33').Value = '  -2.28%  '
$ws.Range('E34').Value = '  -2.35%  '
$ws.Range('D35').Value = '17.49'
$ws.Range('E35').Value = '  +2.07%  '
$ws.Range('D36').Value = '0.0720'
$ws.Range('E36').Value = '  -3.10%  '
$ws.Range('E37').Value = '  +0.37%  '
$ws.Range('D38').Value = '0.115'
$ws.Range('E38').Value = '  -0.39%  '
$ws.Range('D39').Value = '0.103'
$ws.Range('E39').Value = '  -2.95%  '
$ws.Range('D40').Value = '1.79'
$ws.Range('E40').Value = '  -3.29%  '
$ws.Range('D41').Value = '4.06'
$ws.Range('E41').Value = '  -1.77%  '
$ws.Range('D42').Value = '1.940.75'
$ws.Range('E42').Value = '  -2.99%  '
$ws.Range('D43').Value = '0.0281'
$ws.Range('E43').Value = '  -1.83%  '
$ws.Range('D44').Value = '18.58'
$ws.Range('E44').Value = '  -2.83%  '
$ws.Range('D45').Value = '2.17'
$ws.Range('E45').Value = '  -10.38%  '
$ws.Range('D46').Value = '2.89'
$ws.Range('E46').Value = '  -3.57%  '
$ws.Range('D47').Value = '9.64'
$ws.Range('E47').Value = '  -3.62%  '
$ws.Range('D48').Value = '53.81'
$ws.Range('E48').Value = '  +1.32%  '
$ws.Range('D49').Value = '2.470.05'
$ws.Range('E49').Value = '  -1.32%  '
$ws.Range('D50').Value = '71.11'
$ws.Range('E50').Value = '  -1.55%  '
$ws.Range('D51').Value = '91.10'
$ws.Range('E51').Value = '  -1.35%  '
